# Reorder the D/E/F columns (codeforiati:group-name / codeforiati:category-name /
# codeforiati:category-code) on every row, including the header, so that the
# category-code column comes first:
#   before: D=group-name,  E=category-name, F=category-code
#   after:  D=category-code, E=group-name,   F=category-name
# Column G (codeforiati:group-code) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 235

for ($r = 1; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()

    $ws.Cells.Item($r, 4).Value = $f
    $ws.Cells.Item($r, 5).Value = $d
    $ws.Cells.Item($r, 6).Value = $e
}
